$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 63 (Fecha 2022-10-21 / serial 44855),
# which pushes the existing rows 63..74 down to 64..75. Work from the bottom up so
# we never clobber a row before it has been copied to its new destination.
for ($r = 74; $r -ge 63; $r--) {
    $src = $ws.Range("A" + $r + ":R" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":R" + ($r + 1))
    $src.Copy($dst)
}

# Row 63 now holds the new week's record: same market/product/price details as the
# old row 63 (now duplicated at row 64), just a new date.
$ws.Cells.Item(63, 4).Value = 44855
